$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.023.70"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "2.077.81"
$ws.Range("E3").Value = "  +7.00%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("E6").Value = "  -6.26%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.366"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0736"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.96%  "
$ws.Range("E12").Value = "  +4.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.13%  "
$ws.Range("D14").Value = "2.374.06"
$ws.Range("E14").Value = "  +6.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.823"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").Value = "2.075.68"
$ws.Range("E16").Value = "  +6.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.13%  "
$ws.Range("D18").Value = "36.867.26"
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.15%  "
$ws.Range("D20").Value = "0.0₃0816"
$ws.Range("E20").Value = "  -5.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.29%  "
$ws.Range("E30").Value = "  -6.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +18.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0600"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0898"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("E38").Value = "  +12.60%  "
$ws.Range("E39").Value = "  -7.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.49%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0220"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.62%  "
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.83%  "
$ws.Range("E45").Value = "  -4.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0871"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.13%  "
$ws.Range("E47").Value = "  +4.69%  "
$ws.Range("D48").Value = "1.298.49"
$ws.Range("E48").Value = "  -4.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.15%  "
$ws.Range("D50").Value = "2.250.06"
$ws.Range("E50").Value = "  +5.76%  "
$ws.Range("E51").Value = "  -8.50%  "
